$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "wuhan" (sheet1): add 3 days of data (rows 50-52)
# Columns: C=新增病例 D=死亡 E=治愈 F=累计病例 G=累计死亡 H=累计治愈
# ---------------------------------------------------------------------------
$wuhan = $wb.Worksheets.Item("wuhan")

$wuhan.Range("C50").Value = 383
$wuhan.Range("D50").Value = 19
$wuhan.Range("E50").Value = 1535
$wuhan.Range("F50").Value = 47824
$wuhan.Range("G50").Value = 2104
$wuhan.Range("H50").Value = 13328

$wuhan.Range("C51").Value = 313
$wuhan.Range("D51").Value = 28
$wuhan.Range("E51").Value = 2498
$wuhan.Range("F51").Value = 48137
$wuhan.Range("G51").Value = 2132
$wuhan.Range("H51").Value = 15826

$wuhan.Range("C52").Value = 420
$wuhan.Range("D52").Value = 37
$wuhan.Range("E52").Value = 1726
$wuhan.Range("F52").Value = 48557
$wuhan.Range("G52").Value = 2169
$wuhan.Range("H52").Value = 17552

# ---------------------------------------------------------------------------
# Sheet "hubei" (sheet2): add 3 days of data (rows 50-52)
# Columns: C=新增病例 D=死亡 E=治愈 F=新增疑似 G=新增重症(formula)
#          L=累计病例 M=累计死亡 N=累计治愈 O=现有重症 P=危重症 Q=现有住院
# ---------------------------------------------------------------------------
$hubei = $wb.Worksheets.Item("hubei")

$hubei.Range("C50").Value = 409
$hubei.Range("D50").Value = 26
$hubei.Range("E50").Value = 2288
$hubei.Range("F50").Value = 403
$hubei.Range("L50").Value = 65596
$hubei.Range("M50").Value = 2641
$hubei.Range("N50").Value = 23200
$hubei.Range("O50").Value = 6581
$hubei.Range("P50").Value = 1403
$hubei.Range("Q50").Value = 34978
$hubei.Range("G50").Formula = "=(O50+P50)-(O49+P49)"

$hubei.Range("C51").Value = 318
$hubei.Range("D51").Value = 41
$hubei.Range("E51").Value = 3203
$hubei.Range("F51").Value = 332
$hubei.Range("L51").Value = 65914
$hubei.Range("M51").Value = 2682
$hubei.Range("N51").Value = 26403
$hubei.Range("O51").Value = 6270
$hubei.Range("P51").Value = 1363
$hubei.Range("Q51").Value = 32878
$hubei.Range("G51").Formula = "=(O51+P51)-(O50+P50)"

$hubei.Range("C52").Value = 423
$hubei.Range("D52").Value = 45
$hubei.Range("E52").Value = 2492
$hubei.Range("F52").Value = 159
$hubei.Range("L52").Value = 66337
$hubei.Range("M52").Value = 2727
$hubei.Range("N52").Value = 28895
$hubei.Range("O52").Value = 6056
$hubei.Range("P52").Value = 1314
$hubei.Range("Q52").Value = 31064
$hubei.Range("G52").Formula = "=(O52+P52)-(O51+P51)"

# ---------------------------------------------------------------------------
# Sheet "china" (sheet3): fill in previously-missing 累计病例 column (H)
# for rows 47-49 (values had been mis-typed one column to the right), and
# add 3 new days of data (rows 50-52).
# Columns: C=新增病例 D=死亡 E=治愈 F=重症 G=疑似
#          H=累计病例 I=累计死亡 J=累计治愈 K=现有重症 L=共有疑似
# ---------------------------------------------------------------------------
$china = $wb.Worksheets.Item("china")

$china.Range("H47").Value = 77150
$china.Range("I47").Value = 2592
$china.Range("J47").Value = 24734
$china.Range("K47").Value = 9915
$china.Range("L47").Value = 3434

$china.Range("H48").Value = 77658
$china.Range("I48").Value = 2663
$china.Range("J48").Value = 27323
$china.Range("K48").Value = 9126
$china.Range("L48").Value = 2824

$china.Range("H49").Value = 78064
$china.Range("I49").Value = 2715
$china.Range("J49").Value = 29745
$china.Range("K49").Value = 8752
$china.Range("L49").Value = 2491

$china.Range("C50").Value = 433
$china.Range("D50").Value = 29
$china.Range("E50").Value = 2750
$china.Range("F50").Value = -406
$china.Range("G50").Value = 508
$china.Range("H50").Value = 78497
$china.Range("I50").Value = 2744
$china.Range("J50").Value = 32495
$china.Range("K50").Value = 8346
$china.Range("L50").Value = 2358

$china.Range("C51").Value = 327
$china.Range("D51").Value = 44
$china.Range("E51").Value = 3622
$china.Range("F51").Value = -394
$china.Range("G51").Value = 452
$china.Range("H51").Value = 78824
$china.Range("I51").Value = 2788
$china.Range("J51").Value = 36117
$china.Range("K51").Value = 7952
$china.Range("L51").Value = 2308

$china.Range("C52").Value = 427
$china.Range("D52").Value = 47
$china.Range("E52").Value = 2885
$china.Range("F52").Value = -288
$china.Range("G52").Value = 248
$china.Range("H52").Value = 79251
$china.Range("I52").Value = 2835
$china.Range("J52").Value = 39002
$china.Range("K52").Value = 7664
$china.Range("L52").Value = 1418
